# Adds a new "2022-Q3" sheet (fund holdings) right after "总计" and before
# "2021-Q4", and inserts a corresponding summary row on the "总计" sheet.

$wb = $excel.ActiveWorkbook

# Helper: write a value as TEXT (not auto-converted to a number), even
# when it looks numeric (e.g. "013166", "0.38"), without leaving a
# stray NumberFormat-driven style behind on the cell. Routing the
# literal through a formula forces a string result; pasting that back
# as a value keeps the original (unstyled) cell formatting intact.
function Set-TextValue($cell, $val) {
    $cell.Formula = '="' + $val + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" worksheet by copying the existing
#    "2021-Q4" sheet (same column layout/styles) and placing the copy
#    right before it, then renaming + replacing its data.
# ---------------------------------------------------------------------
$quarterSheet = $wb.Worksheets.Item("2021-Q4")
$quarterSheet.Copy($quarterSheet)

$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q3"

# Drop the extra template rows (the source sheet had 9 data rows, the
# new quarter only has 4) -- this also shrinks the sheet's dimension.
$newSheet.Range("A6:H10").Clear()

$fundRows = @(
    @("013166", "东兴宸祥量化混合A", "0.38", "93.87", "1.23", "0.0047", 5),
    @("009327", "东兴兴晟混合A",     "0.38", "79.70", "1.05", "0.0040", 9),
    @("013167", "东兴宸祥量化混合C", "0.08", "93.87", "1.23", "0.0010", 5),
    @("009328", "东兴兴晟混合C",     "0.07", "79.70", "1.05", "0.0007", 9)
)

$r = 2
foreach ($row in $fundRows) {
    Set-TextValue $newSheet.Cells.Item($r, 2) $row[0]
    Set-TextValue $newSheet.Cells.Item($r, 3) $row[1]
    Set-TextValue $newSheet.Cells.Item($r, 4) $row[2]
    Set-TextValue $newSheet.Cells.Item($r, 5) $row[3]
    Set-TextValue $newSheet.Cells.Item($r, 6) $row[4]
    Set-TextValue $newSheet.Cells.Item($r, 7) $row[5]
    $newSheet.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Insert the matching summary row into "总计" (row 2), pushing the
#    existing quarters down and renumbering the index column.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summaryRows = @(
    @(0, "2022-Q3", 4, 0.01),
    @(1, "2021-Q4", 9, 0.97),
    @(2, "2021-Q1", 10, 1.61),
    @(3, "2020-Q4", 2, 0.01)
)

$r = 2
foreach ($row in $summaryRows) {
    $summary.Cells.Item($r, 1).Value = $row[0]
    $summary.Cells.Item($r, 2).Value = $row[1]
    $summary.Cells.Item($r, 3).Value = $row[2]
    $summary.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Row 5 is brand new territory on this sheet -- it didn't inherit the
# "index" column style (A2:A4 = s2) the way the in-range rows did, so
# copy that formatting across explicitly (value is preserved).
$summary.Cells.Item(4, 1).Copy() | Out-Null
$summary.Cells.Item(5, 1).PasteSpecial(-4122) | Out-Null
